$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (interest count) for two events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 9436
$wsExhibit.Range("F5").Value = 521

# Sheet "全部类型" - same two events also appear here, keep in sync
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9436
$wsAll.Range("F5").Value = 521
